# This script updates the "Ticket Sales" (Q) and "Embarking" (R) columns
# for a number of rows in the active worksheet, as part of a debugging
# session (adding/adjusting sample data values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> [new Q value, new R value]
# A value of $null means "leave this cell unchanged"
$updates = @{
    3   = @(72, 30)
    10  = @(88, 55)
    17  = @(97, 54)
    23  = @(40, 27)
    32  = @(80, 1)
    40  = @(94, 19)
    49  = @(18, 18)
    58  = @(69, 4)
    66  = @(34, 2)
    74  = @(52, 31)
    78  = @(22, 15)
    89  = @(44, 31)
    97  = @(17, 4)
    106 = @(82, 37)
    115 = @(95, 20)
    124 = @(40, 1)
    133 = @(86, 5)
    142 = @($null, 54)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $qVal = $vals[0]
    $rVal = $vals[1]

    if ($null -ne $qVal) {
        $ws.Cells.Item($row, 17).Value = $qVal
    }
    if ($null -ne $rVal) {
        $ws.Cells.Item($row, 18).Value = $rVal
    }
}
